$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.299.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.655.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.25%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.21%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  -0.61%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.654.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.28%  '

$ws.Range("E10").Value = '  -2.14%  '

$ws.Range("E11").Value = '  +2.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.356'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.144.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '72.264.61'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.05%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000185'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.651.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '370.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.18'
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.00%  '

$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.797.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.50%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0967'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.38%  '

$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '495.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.53%  '

$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.66%  '

$ws.Range("E38").Value = '  -0.64%  '

$ws.Range("E39").Value = '  -0.38%  '

$ws.Range("E40").Value = '  -2.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.91%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.28%  '

$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '155.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.85%  '

$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("E49").Value = '  +2.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0756'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.33%  '
